$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 612 ("「がんばって！」بالتوفيق" entry) and shift rows below up.
$ws.Rows("612:612").Delete()
